$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I7").Value = "[4]"
$ws.Range("I17").Value = "[12]"
$ws.Range("I25").Value = "[4]"
$ws.Range("I33").Value = "[4]"
$ws.Range("I39").Value = "[12]"
$ws.Range("I41").Value = "[4]"
$ws.Range("I42").Value = "[4]"

$ws.Rows.Item(25).RowHeight = 13.8
$ws.Rows.Item(33).RowHeight = 13.8
$ws.Rows.Item(39).RowHeight = 13.8
$ws.Rows.Item(41).RowHeight = 13.8
$ws.Rows.Item(42).RowHeight = 13.8

$ws.Columns.Item(6).ColumnWidth = 51.2

$ws.Range("I1").Select()
